$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row
$ws.Range("B1").Value = "http://dbpedia.org/ontology/deathPlace"
$ws.Range("C1").Value = "http://dbpedia.org/ontology/parent"
$ws.Range("D1").Value = "http://dbpedia.org/ontology/deathDate"

# Apply the same header style (s=1) used by A1:C1 to the new D1 cell
$ws.Range("A1").Copy()
$ws.Range("D1").PasteSpecial(-4122)  # xlPasteFormats

# Update row 2
$ws.Range("A2").Value = "http://dbpedia.org/resource/Giovanni_Francesco_Guidi_di_Bagno"
$ws.Range("B2").Value = "http://dbpedia.org/resource/Rome"
$ws.Range("C2").Value = "http://dbpedia.org/resource/Colonna_family"
$ws.Range("D2").Value = "http://dbpedia.org/resource/1641"

# Update row 3
$ws.Range("A3").Value = "http://dbpedia.org/resource/Giovanni_Doria"
$ws.Range("B3").Value = "http://dbpedia.org/resource/Palermo"
$ws.Range("C3").Value = "http://dbpedia.org/resource/Giovanni_Andrea_Doria"
$ws.Range("D3").Value = "http://dbpedia.org/resource/1642"

# Remove rows 4-8 (old data no longer present)
$ws.Range("A4:A8").EntireRow.Delete()
